# Finished Week 13 logging
# Update row 2 (Home) stats on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# OFF sheet
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 237
$wsOff.Range("C2").Value = 158
$wsOff.Range("D2").Value = 49
$wsOff.Range("E2").Value = 16
$wsOff.Range("F2").Value = 5

# DEF sheet
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 217
$wsDef.Range("C2").Value = 152
$wsDef.Range("D2").Value = 48
$wsDef.Range("E2").Value = 17
$wsDef.Range("F2").Value = 6
